# Daily refresh of the cryptos price/volume table produced by the scraper
# GitHub Action. Column D ("Price") holds values like "27.402.07" or
# "315.40" that must stay exact strings (trailing zeros, thousands dots,
# etc. matter) - plain `.Value =` would let Excel re-parse them as numbers
# and silently mangle them (e.g. "315.40" -> 315.4), so Set-TextValue
# forces the cell to Text format first, writes the literal string, then
# restores the cell's style to "Normal" so no stray number-format style
# is left behind on cells that were previously unstyled.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2: Bitcoin
Set-TextValue "D2" "27.402.07"
$ws.Range("E2").Value = "  +1.30%  "
# Row 3: Ethereum
Set-TextValue "D3" "1.862.59"
$ws.Range("E3").Value = "  +1.84%  "
# Row 4: TetherUSD
Set-TextValue "D4" "1.002"
$ws.Range("E4").Value = "  -0.59%  "
# Row 5: BNB
Set-TextValue "D5" "315.40"
$ws.Range("E5").Value = "  +0.82%  "
# Row 6: USDC
Set-TextValue "D6" "1.002"
$ws.Range("E6").Value = "  -0.55%  "
# Row 7: XRP
Set-TextValue "D7" "0.4627"
$ws.Range("E7").Value = "  +0.47%  "
# Row 8: Cardano
Set-TextValue "D8" "0.3717"
$ws.Range("E8").Value = "  +0.56%  "
# Row 9: Dogecoin
Set-TextValue "D9" "0.07323"
$ws.Range("E9").Value = "  -0.38%  "
# Row 10: Polygon
Set-TextValue "D10" "0.8900"
$ws.Range("E10").Value = "  +2.11%  "
# Row 11: Solana
Set-TextValue "D11" "20.09"
$ws.Range("E11").Value = "  +1.51%  "
# Row 12: TRON
Set-TextValue "D12" "0.07838"
$ws.Range("E12").Value = "  -1.42%  "
# Rows 13-16: ranking reshuffled - WrappedEther jumps to #13 (was #16),
# pushing Polkadot/Chainlink/Litecoin each down one spot.
# Row 13: was Polkadot -> now WrappedEther
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D13" "1.835.17"
$ws.Range("E13").Value = "  +3.08%  "
# Row 14: was Chainlink -> now Polkadot
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D14" "5.402"
$ws.Range("E14").Value = "  +1.16%  "
# Row 15: was Litecoin -> now Chainlink
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D15" "6.558"
$ws.Range("E15").Value = "  +0.36%  "
# Row 16: was WrappedEther -> now Litecoin
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D16" "91.87"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("E17").Value = "  -0.52%  "
# Row 18: ShibaInu
Set-TextValue "D18" "0.000008981"
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("E19").Value = "  -0.26%  "
# Row 20: Avalanche
Set-TextValue "D20" "14.82"
$ws.Range("E20").Value = "  +0.94%  "
# Row 21: WrappedBTC
Set-TextValue "D21" "27.409.54"
$ws.Range("E21").Value = "  +2.66%  "
# Row 22: Uniswap
Set-TextValue "D22" "5.136"
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("E23").Value = "  -0.10%  "
# Row 24: WrappedliquidstakedEther2.0
Set-TextValue "D24" "2.032.78"
$ws.Range("E24").Value = "  +5.44%  "
# Row 25: Toncoin
Set-TextValue "D25" "1.937"
$ws.Range("E25").Value = "  +4.58%  "
# Row 26: Monero
Set-TextValue "D26" "152.06"
$ws.Range("E26").Value = "  -0.27%  "
# Row 27: EthereumClassic
Set-TextValue "D27" "18.47"
$ws.Range("E27").Value = "  -0.40%  "
# Row 28: LidoDAOToken
Set-TextValue "D28" "2.057"
$ws.Range("E28").Value = "  -0.70%  "
# Row 29: InternetComputer(DFINITY)
Set-TextValue "D29" "5.108"
$ws.Range("E29").Value = "  +0.25%  "
# Row 30: BitcoinCash
Set-TextValue "D30" "116.23"
$ws.Range("E30").Value = "  +0.81%  "
# Row 31: Stellar
Set-TextValue "D31" "0.08850"
$ws.Range("E31").Value = "  -0.30%  "
# Row 32: HuobiToken
Set-TextValue "D32" "3.131"
$ws.Range("E32").Value = "  +5.17%  "
# Row 33: ImmutableX
Set-TextValue "D33" "0.7717"
$ws.Range("E33").Value = "  +5.51%  "
# Row 34: ARBITRUM
Set-TextValue "D34" "1.175"
$ws.Range("E34").Value = "  +3.42%  "
# Rows 35-36: RenderToken/Filecoin swap ranking positions
# Row 35: was Filecoin -> now RenderToken
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D35" "2.799"
$ws.Range("E35").Value = "  +13.88%  "
# Row 36: was RenderToken -> now Filecoin
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D36" "4.519"
$ws.Range("E36").Value = "  +1.74%  "
# Row 37: TrustWalletToken
Set-TextValue "D37" "1.080"
$ws.Range("E37").Value = "  +0.86%  "
# Row 38: VeChain
Set-TextValue "D38" "0.01961"
$ws.Range("E38").Value = "  +0.91%  "
$ws.Range("E39").Value = "  +0.08%  "
# Row 40: MXToken
Set-TextValue "D40" "2.980"
$ws.Range("E40").Value = "  +1.30%  "
# Row 41: FraxShare
Set-TextValue "D41" "7.086"
$ws.Range("E41").Value = "  -0.45%  "
# Row 42: TheSandbox
Set-TextValue "D42" "0.5151"
$ws.Range("E42").Value = "  -0.37%  "
# Row 43: Algorand
Set-TextValue "D43" "0.1642"
$ws.Range("E43").Value = "  +0.71%  "
# Row 44: Aptos
Set-TextValue "D44" "8.419"
$ws.Range("E44").Value = "  +2.48%  "
# Row 45: Decentraland
Set-TextValue "D45" "0.4816"
$ws.Range("E45").Value = "  -0.29%  "
# Row 46: EnergySwap
Set-TextValue "D46" "10.35"
$ws.Range("E46").Value = "  +0.82%  "
$ws.Range("E47").Value = "  -0.55%  "
# Row 48: Quant
Set-TextValue "D48" "103.14"
$ws.Range("E48").Value = "  +1.06%  "
$ws.Range("E49").Value = "  +1.67%  "
# Row 50: Cronos
Set-TextValue "D50" "0.06220"
$ws.Range("E50").Value = "  -0.13%  "
# Row 51: Aave
Set-TextValue "D51" "65.52"
$ws.Range("E51").Value = "  +1.13%  "
